$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-06-11T08:08:31+00:00"

# --- Concepts sheet: fix accents / typos in Display column ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C2").Value  = "spécialité"
$concepts.Range("C5").Value  = "préparation hospitalière"
$concepts.Range("C6").Value  = "importation"
$concepts.Range("C8").Value  = "préparation magistrale"
$concepts.Range("C9").Value  = "matière première"
$concepts.Range("C10").Value = "médicament virtuel THESORIMED"
$concepts.Range("C11").Value = "médicament virtuel THERIAQUE"
$concepts.Range("C12").Value = "médicament virtuel BCB DEXTER"
$concepts.Range("C13").Value = "médicament virtuel VIDAL"
$concepts.Range("C14").Value = "dénomination commune"
$concepts.Range("C15").Value = "médicament virtuel MedicaBase"
